# feat: add 2022-Q1 data
#
# 1. Insert a new top data row into the "总计" (summary) sheet for 2022-Q1,
#    shifting the existing rows down.
# 2. Insert a new worksheet "2022-Q1" (fund-holding detail) positioned
#    between "2021-Q4" and "总计", and populate it with the same column
#    layout as the other quarter sheets.
#
# NOTE: worksheet variables in this host resolve by *position*, not by a
# stable object identity. Any sheet-insertion that changes where "总计"
# sits in the tab order would silently redirect an already-captured
# `$wsTotal` reference to whatever sheet now occupies that slot. To avoid
# that, all edits to "总计" happen first, while it is still the last
# sheet, before the new sheet is created in front of it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the 2022-Q1 summary row at the top of "总计"
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

$wsTotal.Rows.Item(2).Insert(-4121)
$wsTotal.Range("B2:D2").ClearFormats()

$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)
$wsTotal.Application.CutCopyMode = $false

$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Cells.Item(2, 2).Value = "2022-Q1"
$wsTotal.Cells.Item(2, 3).Value = 3
$wsTotal.Cells.Item(2, 4).Value = 0.8

# Column A is a plain 0-based row counter, independent of the row's data;
# the rows that were pushed down must have their index renumbered too.
$wsTotal.Cells.Item(3, 1).Value = 1
$wsTotal.Cells.Item(4, 1).Value = 2

$wsTotal.Range("A1").Select()

# ---------------------------------------------------------------------
# 2. Create & populate the "2022-Q1" sheet right after "2021-Q4"
# ---------------------------------------------------------------------
$wsQ4 = $wb.Worksheets.Item("2021-Q4")
$wsNew = $wb.Worksheets.Add($null, $wsQ4)
$wsNew.Name = "2022-Q1"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 2 + $i
    $wsNew.Cells.Item(1, $col).Value = $headers[$i]
}

# Copy the header look & feel (bold font, thin border, centered) from the
# "2021-Q4" sheet so the new sheet matches the workbook's existing style.
$wsQ4 = $wb.Worksheets.Item("2021-Q4")
$wsQ4.Range("B1:H1").Copy()
$wsNew.Range("B1:H1").PasteSpecial(-4122)
$wsNew.Application.CutCopyMode = $false

function Set-TextValue($range, [string]$text) {
    # Force the value to be stored as text, even when it looks numeric,
    # matching the source data (e.g. "8.67", "012850"), then drop the
    # temporary "@" number format so no style residue is left behind.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.NumberFormat = "General"
    $range.ClearFormats()
}

$rows = @(
    @{ B = "012850"; C = "中融低碳经济3个月持有期混合型证券投资基金A"; D = "8.67";  E = "65.64"; F = "4.50"; G = "0.3902"; H = 2 },
    @{ B = "009011"; C = "华夏睿阳一年持有期混合";                      D = "17.06"; E = "82.70"; F = "2.00"; G = "0.3412"; H = 9 },
    @{ B = "012851"; C = "中融低碳经济3个月持有期混合型证券投资基金C"; D = "1.49";  E = "65.64"; F = "4.50"; G = "0.0670"; H = 2 }
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    # Column A: the 0-based row index, styled like the header (bold/border).
    $wsQ4 = $wb.Worksheets.Item("2021-Q4")
    $wsQ4.Range("A2").Copy()
    $wsNew.Range("A$r").PasteSpecial(-4122)
    $wsNew.Application.CutCopyMode = $false
    $wsNew.Cells.Item($r, 1).Value = $i

    Set-TextValue $wsNew.Range("B$r") $data.B
    Set-TextValue $wsNew.Range("C$r") $data.C
    Set-TextValue $wsNew.Range("D$r") $data.D
    Set-TextValue $wsNew.Range("E$r") $data.E
    Set-TextValue $wsNew.Range("F$r") $data.F
    Set-TextValue $wsNew.Range("G$r") $data.G

    $wsNew.Cells.Item($r, 8).Value = $data.H
}

$wsNew.Range("A1").Select()
